$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 250
$ws.Range("I11").Value = 250
$ws.Range("K11").Value = 250
$ws.Range("M11").Value = -110

$ws.Range("H17").Value = 1517217.5
$ws.Range("J17").Value = 1564609.2
$ws.Range("L17").Value = 4693827.6
$ws.Range("N17").Value = -4694163.6

$ws.Range("H76").Value = 2908.8
$ws.Range("I76").Value = 2908.8
$ws.Range("K76").Value = 2908.8
$ws.Range("M76").Value = -2593.8

$ws.Range("H79").Value = 2908.8
$ws.Range("I79").Value = 2908.8
$ws.Range("K79").Value = 2908.8
$ws.Range("M79").Value = -1816.8

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = ""
$ws.Range("N115").Value = ""

$ws.Range("H132").Value = 1834.0588
$ws.Range("I132").Value = 1830.74
$ws.Range("K132").Value = 5492.22
$ws.Range("M132").Value = -2962.22

$ws.Range("H138").Value = 1572.8025
$ws.Range("I138").Value = 638.4
$ws.Range("J138").Value = 2283.761
$ws.Range("K138").Value = 1915.2
$ws.Range("L138").Value = 6851.282999999999
$ws.Range("M138").Value = 3224.8
$ws.Range("N138").Value = -17131.283

$ws.Range("H141").Value = 2026.6666
$ws.Range("I141").Value = 1780
$ws.Range("K141").Value = 5340
$ws.Range("M141").Value = -160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4132.433
$ws.Range("I32").Value = 4177.357
$ws.Range("K32").Value = 4177.357
$ws.Range("M32").Value = -3890.357

$ws.Range("H132").Value = 11338.039
$ws.Range("I132").Value = 1350.6364
$ws.Range("J132").Value = 74116
$ws.Range("K132").Value = 4051.9092
$ws.Range("L132").Value = 222348
$ws.Range("M132").Value = -1521.9092
$ws.Range("N132").Value = -227408

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J107").Value = 1200
$ws.Range("L107").Value = 1200
$ws.Range("N107").Value = -5040

$ws.Range("H134").Value = 4633.64
$ws.Range("I134").Value = 4897.4346
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 14692.3038
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -12157.3038
$ws.Range("N134").Value = -9870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 49.5
$ws.Range("I7").Value = 57.333332
$ws.Range("J7").Value = 35.4
$ws.Range("K7").Value = 57.333332
$ws.Range("L7").Value = 35.4
$ws.Range("M7").Value = 55.666668
$ws.Range("N7").Value = -261.4

$ws.Range("H31").Value = 13605.132
$ws.Range("I31").Value = 23811.166
$ws.Range("J31").Value = 4419.7
$ws.Range("K31").Value = 23811.166
$ws.Range("L31").Value = 4419.7
$ws.Range("M31").Value = -23516.166
$ws.Range("N31").Value = -5009.7

$ws.Range("H34").Value = 13605.132
$ws.Range("I34").Value = 23811.166
$ws.Range("J34").Value = 4419.7
$ws.Range("K34").Value = 23811.166
$ws.Range("L34").Value = 4419.7
$ws.Range("M34").Value = -23609.166
$ws.Range("N34").Value = -4823.7

$ws.Range("H86").Value = 8343376.5
$ws.Range("I86").Value = 3099.9167
$ws.Range("J86").Value = 20853792
$ws.Range("K86").Value = 3099.9167
$ws.Range("L86").Value = 20853792
$ws.Range("M86").Value = -1976.9167
$ws.Range("N86").Value = -20856038

$ws.Range("H89").Value = 8343376.5
$ws.Range("I89").Value = 3099.9167
$ws.Range("J89").Value = 20853792
$ws.Range("K89").Value = 15499.5835
$ws.Range("L89").Value = 104268960
$ws.Range("M89").Value = -9883.583500000001
$ws.Range("N89").Value = -104280192

$ws.Range("H94").Value = 2925.8262
$ws.Range("I94").Value = 1605.7693
$ws.Range("K94").Value = 1605.7693
$ws.Range("M94").Value = -1154.7693

$ws.Range("H122").Value = 927.2174
$ws.Range("I122").Value = 947.7273
$ws.Range("K122").Value = 2843.1819
$ws.Range("M122").Value = -393.1819

$ws.Range("H132").Value = 11087.036
$ws.Range("I132").Value = 13623.122
$ws.Range("J132").Value = 3659.9285
$ws.Range("K132").Value = 40869.36599999999
$ws.Range("L132").Value = 10979.7855
$ws.Range("M132").Value = -38339.36599999999
$ws.Range("N132").Value = -16039.7855

$ws.Range("H134").Value = 1045.8334
$ws.Range("I134").Value = 769.3095
$ws.Range("K134").Value = 2307.9285
$ws.Range("M134").Value = 227.0715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 589.1667
$ws.Range("J60").Value = 1500
$ws.Range("L60").Value = 4500
$ws.Range("N60").Value = -5002

$ws.Range("H98").Value = 1349
$ws.Range("I98").Value = 1748.5
$ws.Range("J98").Value = 949.5
$ws.Range("K98").Value = 5245.5
$ws.Range("L98").Value = 2848.5
$ws.Range("M98").Value = -3747.5
$ws.Range("N98").Value = -5844.5

$ws.Range("H122").Value = 1503.579
$ws.Range("J122").Value = 1635.1765
$ws.Range("L122").Value = 14716.5885
$ws.Range("N122").Value = -19616.5885

$ws.Range("H131").Value = 102836.055
$ws.Range("I131").Value = 666.5
$ws.Range("J131").Value = 104964.586
$ws.Range("K131").Value = 1999.5
$ws.Range("L131").Value = 314893.758
$ws.Range("M131").Value = 3040.5
$ws.Range("N131").Value = -324973.758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6334.0835
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 5668.1665
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 5668.1665
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -6208.1665

$ws.Range("H73").Value = 6334.0835
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 5668.1665
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 5668.1665
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -7540.1665

$ws.Range("H122").Value = 51283070
$ws.Range("I122").Value = 22223476
$ws.Range("J122").Value = 90909790
$ws.Range("K122").Value = 66670428
$ws.Range("L122").Value = 272729370
$ws.Range("M122").Value = -66667978
$ws.Range("N122").Value = -272734270

$ws.Range("H132").Value = 17550.828
$ws.Range("I132").Value = 3160.4285
$ws.Range("K132").Value = 9481.2855
$ws.Range("M132").Value = -6951.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2415.5
$ws.Range("I22").Value = 2376.9092
$ws.Range("J22").Value = 2500.4
$ws.Range("K22").Value = 2376.9092
$ws.Range("L22").Value = 2500.4
$ws.Range("M22").Value = -2081.9092
$ws.Range("N22").Value = -3090.4

$ws.Range("H27").Value = 2415.5
$ws.Range("I27").Value = 2376.9092
$ws.Range("J27").Value = 2500.4
$ws.Range("K27").Value = 2376.9092
$ws.Range("L27").Value = 2500.4
$ws.Range("M27").Value = -2269.9092
$ws.Range("N27").Value = -2714.4

$ws.Range("H46").Value = 859.26666
$ws.Range("I46").Value = 765.44446
$ws.Range("K46").Value = 765.44446
$ws.Range("M46").Value = -577.44446

$ws.Range("H55").Value = 53.35294
$ws.Range("I55").Value = 45.444443
$ws.Range("K55").Value = 45.444443
$ws.Range("M55").Value = 127.555557

$ws.Range("H132").Value = 2013.7556
$ws.Range("I132").Value = 923.13794
$ws.Range("J132").Value = 3990.5
$ws.Range("K132").Value = 2769.41382
$ws.Range("L132").Value = 11971.5
$ws.Range("M132").Value = -239.4138199999998
$ws.Range("N132").Value = -17031.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1229182
$ws.Range("I113").Value = 960.7692
$ws.Range("J113").Value = 3003279.2
$ws.Range("K113").Value = 2882.3076
$ws.Range("L113").Value = 9009837.600000001
$ws.Range("M113").Value = -712.3076000000001
$ws.Range("N113").Value = -9014177.600000001

$ws.Range("H122").Value = 1709.7354
$ws.Range("I122").Value = 1717.8
$ws.Range("J122").Value = 1687.3334
$ws.Range("K122").Value = 5153.4
$ws.Range("L122").Value = 5062.0002
$ws.Range("M122").Value = -2703.4
$ws.Range("N122").Value = -9962.0002
